$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_edited_time (D7:D12)
$ws.Range("D7").Value = "2024-07-21T16:44:00.000Z"
$ws.Range("D8").Value = "2024-07-21T16:44:00.000Z"
$ws.Range("D9").Value = "2024-07-21T16:44:00.000Z"
$ws.Range("D10").Value = "2024-07-21T16:44:00.000Z"
$ws.Range("D11").Value = "2024-07-21T16:44:00.000Z"
$ws.Range("D12").Value = "2024-07-21T16:44:00.000Z"

# Update numeric values in row 7
$ws.Range("W7").Value = 225658000
$ws.Range("AA7").Value = 215160000
$ws.Range("AE7").Value = 440818000
$ws.Range("AH7").Value = 376818000
$ws.Range("AK7").Value = 61
$ws.Range("AQ7").Value = 414118000
